$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "60.895.69"
$ws.Range("E2").Value = "  -0.84%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.389.85"
$ws.Range("E3").Value = "  -1.31%  "

# Row 4 - TetherUSD (numeric-looking, force text)
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.01%  "

# Row 5 - BNB (price unchanged)
$ws.Range("E5").Value = "  -0.73%  "

# Row 6 - Solana (numeric-looking, force text)
$ws.Range("D6").Value = "'141.87"
$ws.Range("E6").Value = "  -2.20%  "

# Row 7 - LidoStakedEther
$ws.Range("D7").Value = "3.390.45"
$ws.Range("E7").Value = "  -1.29%  "

# Row 8 - USDC (price unchanged)
$ws.Range("E8").Value = "  +0.03%  "

# Row 9 - XRP (price unchanged)
$ws.Range("E9").Value = "  -0.16%  "

# Row 10 - Toncoin (numeric-looking, force text)
$ws.Range("D10").Value = "'7.54"
$ws.Range("E10").Value = "  -1.70%  "

# Row 11 - Dogecoin (price unchanged)
$ws.Range("E11").Value = "  -1.79%  "

# Row 12 - Cardano (price unchanged)
$ws.Range("E12").Value = "  +2.28%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "3.968.44"
$ws.Range("E13").Value = "  -1.30%  "

# Row 14 - TRON (price unchanged)
$ws.Range("E14").Value = "  +2.07%  "

# Row 15 - Avalanche (price unchanged)
$ws.Range("E15").Value = "  +0.23%  "

# Row 16 - ShibaInu (price unchanged)
$ws.Range("E16").Value = "  -1.03%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "3.395.16"
$ws.Range("E17").Value = "  -1.13%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "60.968.67"
$ws.Range("E18").Value = "  -0.91%  "

# Row 19 - Polkadot (numeric-looking, force text)
$ws.Range("D19").Value = "'6.16"
$ws.Range("E19").Value = "  -2.11%  "

# Row 20 - Chainlink (price unchanged)
$ws.Range("E20").Value = "  -2.59%  "

# Row 21 - Uniswap (numeric-looking, force text)
$ws.Range("D21").Value = "'8.97"
$ws.Range("E21").Value = "  -4.58%  "

# Row 22 - BitcoinCash (numeric-looking, force text)
$ws.Range("D22").Value = "'384.16"
$ws.Range("E22").Value = "  -2.99%  "

# Row 23 - Polygon (price unchanged)
$ws.Range("E23").Value = "  -1.50%  "

# Row 24 - Litecoin (numeric-looking, force text)
$ws.Range("D24").Value = "'74.45"
$ws.Range("E24").Value = "  +1.00%  "

# Row 25 - Dai (price unchanged)
$ws.Range("E25").Value = "  +0.46%  "

# Row 26 - PEPE (price unchanged)
$ws.Range("E26").Value = "  -4.59%  "

# Row 27 - WrappedeETH
$ws.Range("D27").Value = "3.527.70"

# Row 28 - Kaspa (price unchanged)
$ws.Range("E28").Value = "  -1.36%  "

# Row 29 - Binance-PegBSC-USD (price unchanged)
$ws.Range("E29").Value = "  -0.25%  "

# Row 30 - RenderToken (price unchanged)
$ws.Range("E30").Value = "  -2.75%  "

# Row 31 - InternetComputer(DFINITY) (numeric-looking, force text)
$ws.Range("D31").Value = "'7.97"
$ws.Range("E31").Value = "  -3.38%  "

# Row 33 - Fetch.AI (numeric-looking, force text)
$ws.Range("D33").Value = "'1.43"
$ws.Range("E33").Value = "  -2.74%  "

# Row 34 - USDe (price unchanged)
$ws.Range("E34").Value = "  +0.00%  "

# Row 35 - EthereumClassic (numeric-looking, force text)
$ws.Range("D35").Value = "'23.53"
$ws.Range("E35").Value = "  -1.71%  "

# Row 36 - Aptos (price unchanged)
$ws.Range("E36").Value = "  -0.34%  "

# Row 37 - Monero (numeric-looking, force text)
$ws.Range("D37").Value = "'167.58"
$ws.Range("E37").Value = "  +0.06%  "

# Row 38 - RenzoRestakedETH
$ws.Range("D38").Value = "3.420.23"
$ws.Range("E38").Value = "  -1.25%  "

# Row 39 - NEARProtocol (numeric-looking, force text)
$ws.Range("D39").Value = "'4.98"
$ws.Range("E39").Value = "  -2.42%  "

# Row 40 - ImmutableX (price unchanged)
$ws.Range("E40").Value = "  -4.21%  "

# Row 41 - Hedera (numeric-looking, force text)
$ws.Range("D41").Value = "'0.0775"
$ws.Range("E41").Value = "  -1.26%  "

# Row 42 - EnergySwap (numeric-looking, force text)
$ws.Range("D42").Value = "'27.30"
$ws.Range("E42").Value = "  +0.99%  "

# Row 43 - FirstDigitalUSD (price unchanged)
$ws.Range("E43").Value = "  -0.03%  "

# Row 44 - Mantle (price unchanged)
$ws.Range("E44").Value = "  -2.46%  "

# Row 45 - OKB (numeric-looking, force text)
$ws.Range("D45").Value = "'42.16"
$ws.Range("E45").Value = "  -0.25%  "

# Row 46 - Filecoin (price unchanged)
$ws.Range("E46").Value = "  -1.27%  "

# Row 47 - Stacks (price unchanged)
$ws.Range("E47").Value = "  -3.39%  "

# Row 48 - ONDO (price unchanged)
$ws.Range("E48").Value = "  -1.13%  "

# Row 49 - Maker
$ws.Range("D49").Value = "2.476.60"
$ws.Range("E49").Value = "  -4.67%  "

# Row 50 - Cosmos (price unchanged)
$ws.Range("E50").Value = "  -1.49%  "

# Row 51 - InjectiveProtocol (numeric-looking, force text)
$ws.Range("D51").Value = "'23.01"
$ws.Range("E51").Value = "  -1.04%  "
